$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.682.04'
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').Value = '1.867.78'
$ws.Range('E3').Value = '  -0.70%  '
$ws.Range('D4').Value = '1.011'
$ws.Range('E4').Value = '  +0.48%  '
$ws.Range('D5').Value = '333.89'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '0.3923'
$ws.Range('E8').Value = '  -0.30%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = '0.07999'
$ws.Range('E9').Value = '  -0.89%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').Value = '45.12'
$ws.Range('E10').Value = '  -4.94%  '
$ws.Range('D11').Value = '1.004'
$ws.Range('E11').Value = '  -2.33%  '
$ws.Range('D12').Value = '21.88'
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('D13').Value = '6.006'
$ws.Range('D14').Value = '1.866.34'
$ws.Range('E14').Value = '  -0.97%  '
$ws.Range('D15').Value = '7.265'
$ws.Range('E15').Value = '  +1.70%  '
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').Value = '88.66'
$ws.Range('E17').Value = '  +1.76%  '
$ws.Range('D18').Value = '0.06743'
$ws.Range('E18').Value = '  +0.56%  '
$ws.Range('D19').Value = '0.00001044'
$ws.Range('E19').Value = '  -0.43%  '
$ws.Range('D20').Value = '17.16'
$ws.Range('E20').Value = '  -1.10%  '
$ws.Range('E21').Value = '  +0.39%  '
$ws.Range('D22').Value = '27.660.99'
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('D23').Value = '5.472'
$ws.Range('E23').Value = '  -1.31%  '
$ws.Range('E24').Value = '  -1.01%  '
$ws.Range('D25').Value = '2.314'
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').Value = '2.084.16'
$ws.Range('E26').Value = '  -1.15%  '
$ws.Range('D27').Value = '159.38'
$ws.Range('E27').Value = '  -0.40%  '
$ws.Range('D28').Value = '19.83'
$ws.Range('E28').Value = '  -1.75%  '
$ws.Range('D29').Value = '2.161'
$ws.Range('E29').Value = '  +2.58%  '
$ws.Range('D30').Value = '5.459'
$ws.Range('E30').Value = '  -2.35%  '
$ws.Range('D31').Value = '121.92'
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('D32').Value = '0.9842'
$ws.Range('E32').Value = '  -0.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09510'
$ws.Range('E33').Value = '  +0.27%  '
$ws.Range('D34').Value = '3.611'
$ws.Range('E34').Value = '  -0.34%  '
$ws.Range('D35').Value = '5.321'
$ws.Range('E35').Value = '  -0.76%  '
$ws.Range('D36').Value = '1.342'
$ws.Range('E36').Value = '  -7.68%  '
$ws.Range('D37').Value = '0.06065'
$ws.Range('E37').Value = '  -1.20%  '
$ws.Range('D38').Value = '0.02237'
$ws.Range('E38').Value = '  -1.31%  '
$ws.Range('D39').Value = '8.349'
$ws.Range('E39').Value = '  +2.54%  '
$ws.Range('D40').Value = '1.197'
$ws.Range('E41').Value = '  +0.29%  '
$ws.Range('D42').Value = '0.5986'
$ws.Range('E42').Value = '  -0.33%  '
$ws.Range('D43').Value = '0.1889'
$ws.Range('E43').Value = '  -0.66%  '
$ws.Range('D44').Value = '10.31'
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('D45').Value = '1.247'
$ws.Range('E45').Value = '  -0.81%  '
$ws.Range('D46').Value = '0.5659'
$ws.Range('D47').Value = '12.25'
$ws.Range('E47').Value = '  -0.10%  '
$ws.Range('D48').Value = '1.926'
$ws.Range('E48').Value = '  -1.09%  '
$ws.Range('D49').Value = '0.06758'
$ws.Range('E49').Value = '  -2.11%  '
$ws.Range('D50').Value = '112.08'
$ws.Range('E50').Value = '  -2.02%  '
$ws.Range('D51').Value = '3.052'
$ws.Range('E51').Value = '  -10.21%  '
